$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new note text in cell C9 (new shared string entry)
$ws.Range("C9").Value = "C,D không dẫn ra được F"

# Update the active selection to match the target state (F6)
$ws.Range("F6").Select()
